# Updates cryptos list figures (price & 1h volume change) for 22 of 50
# coins, and fixes the ordering / re-maps of 6 coins whose rank changed
# (rows 32/33, 41/42, 49/50) plus a rank-50 swap-in (CoreDAO replacing
# Monero) on row 51, matching an automated "Updated cryptos list" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are quote-prefixed ('...) so Excel keeps them
# as text instead of auto-coercing numeric-looking strings (matching
# prices like "64.520.08" that use "." as a thousands separator and
# would otherwise be misparsed or have trailing zeros stripped).

$ws.Range("D2").Value = "'64.520.08"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "'3.168.80"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'599.08"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "'151.50"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'3.162.95"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("E12").Value = "  -6.13%  "
$ws.Range("E13").Value = "  -5.50%  "
$ws.Range("D14").Value = "'36.86"
$ws.Range("E14").Value = "  -5.30%  "
$ws.Range("D15").Value = "'3.700.67"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "'64.507.48"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").Value = "'3.174.84"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").Value = "'479.64"
$ws.Range("E20").Value = "  -6.03%  "
$ws.Range("D21").Value = "'14.78"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "'0.714"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").Value = "'7.73"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "'13.87"
$ws.Range("E24").Value = "  -5.58%  "
$ws.Range("D25").Value = "'84.37"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("E28").Value = "  -5.98%  "
$ws.Range("E29").Value = "  -5.67%  "
$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  +14.71%  "
$ws.Range("D31").Value = "'7.00"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.72"
$ws.Range("E33").Value = "  -8.22%  "
$ws.Range("D34").Value = "'26.84"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("D36").Value = "'6.11"
$ws.Range("E36").Value = "  -5.96%  "
$ws.Range("D37").Value = "'54.42"
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("D39").Value = "'0.0₃0737"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("D40").Value = "'456.87"
$ws.Range("E40").Value = "  -10.09%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0402"
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.124"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").Value = "'8.49"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").Value = "'2.868.21"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  -8.20%  "
$ws.Range("D47").Value = "'27.22"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.116"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.51"
$ws.Range("E51").Value = "  +1.04%  "
